# Weekly data refresh: insert this week's two new "Cebollín" price rows
# (Primera / Segunda) at the top of the Vega Monumental Concepción block,
# pushing the existing rows 81-99 down to 83-101.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 81:82, shifting rows 81-99 down to 83-101.
$ws.Range("A81:A82").EntireRow.Insert()

# New row 81 - "Primera" quality, week of 45007.
$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(81, 3).Value = "Bíobío"
$ws.Cells.Item(81, 4).Value = 45007
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 100112037
$ws.Cells.Item(81, 7).Value = "Cebollín"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 200
$ws.Cells.Item(81, 11).Value = 700
$ws.Cells.Item(81, 12).Value = 800
$ws.Cells.Item(81, 13).Value = 750
$ws.Cells.Item(81, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(81, 15).Value = "Región de Ñuble"
$ws.Cells.Item(81, 16).Value = 125
$ws.Cells.Item(81, 17).Value = 6
$ws.Cells.Item(81, 18).Value = "Hortaliza"

# New row 82 - "Segunda" quality, week of 45007.
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(82, 3).Value = "Bíobío"
$ws.Cells.Item(82, 4).Value = 45007
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = 100112037
$ws.Cells.Item(82, 7).Value = "Cebollín"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Segunda"
$ws.Cells.Item(82, 10).Value = 100
$ws.Cells.Item(82, 11).Value = 600
$ws.Cells.Item(82, 12).Value = 600
$ws.Cells.Item(82, 13).Value = 600
$ws.Cells.Item(82, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(82, 15).Value = "Región de Ñuble"
$ws.Cells.Item(82, 16).Value = 100
$ws.Cells.Item(82, 17).Value = 6
$ws.Cells.Item(82, 18).Value = "Hortaliza"
